$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update IPC PO (C), DELTA (D) and DELTA^2 (E) for each data row (2-51)
$ws.Range("C2").Value = 28.56809871931244
$ws.Range("D2").Value = -0.7819012806875634
$ws.Range("E2").Value = 0.6113696127408518
$ws.Range("C3").Value = 28.56090230188462
$ws.Range("D3").Value = -0.8090976981153766
$ws.Range("E3").Value = 0.6546390850956011
$ws.Range("C4").Value = 29.40648462030243
$ws.Range("D4").Value = -0.1335153796975703
$ws.Range("E4").Value = 0.01782635661578637
$ws.Range("C5").Value = 30.00421754041881
$ws.Range("D5").Value = 0.4542175404188136
$ws.Range("E5").Value = 0.2063135740241165
$ws.Range("C6").Value = 30.08550511156546
$ws.Range("D6").Value = 0.3355051115654604
$ws.Range("E6").Value = 0.112563679886552
$ws.Range("C7").Value = 30.13760866815252
$ws.Range("D7").Value = 0.2976086681525238
$ws.Range("E7").Value = 0.08857091935951907
$ws.Range("C8").Value = 30.24145396543597
$ws.Range("D8").Value = 0.4314539654359741
$ws.Range("E8").Value = 0.1861525242904268
$ws.Range("C9").Value = 30.23881880015279
$ws.Range("D9").Value = 0.3188188001527905
$ws.Range("E9").Value = 0.101645427330865
$ws.Range("C10").Value = 30.3324748470579
$ws.Range("D10").Value = 0.3524748470578949
$ws.Range("E10").Value = 0.1242385178084864
$ws.Range("C11").Value = 30.15780981532706
$ws.Range("D11").Value = 0.1178098153270639
$ws.Range("E11").Value = 0.0138791525873969
$ws.Range("C12").Value = 30.00645248447812
$ws.Range("D12").Value = -0.2035475155218762
$ws.Range("E12").Value = 0.04143159107512844
$ws.Range("C13").Value = 30.11347628287443
$ws.Range("D13").Value = -0.1065237171255689
$ws.Range("E13").Value = 0.01134730231024821
$ws.Range("C14").Value = 29.99755940397955
$ws.Range("D14").Value = -0.382440596020448
$ws.Range("E14").Value = 0.1462608094844755
$ws.Range("C15").Value = 30.02760648670544
$ws.Range("D15").Value = -0.412393513294564
$ws.Range("E15").Value = 0.1700684098074337
$ws.Range("C16").Value = 30.198318076291
$ws.Range("D16").Value = -0.2816819237090051
$ws.Range("E16").Value = 0.07934470614440579
$ws.Range("C17").Value = 30.46317541374856
$ws.Range("D17").Value = -0.2268245862514426
$ws.Range("E17").Value = 0.05144939292813814
$ws.Range("C18").Value = 30.7874844356307
$ws.Range("D18").Value = 0.03748443563070225
$ws.Range("E18").Value = 0.001405082914552261
$ws.Range("C19").Value = 31.48652885860169
$ws.Range("D19").Value = 0.5465288586016932
$ws.Range("E19").Value = 0.2986937932844695
$ws.Range("C20").Value = 31.39704363937706
$ws.Range("D20").Value = 0.4470436393770605
$ws.Range("E20").Value = 0.1998480155074873
$ws.Range("C21").Value = 31.22176637042529
$ws.Range("D21").Value = 0.2017663704252861
$ws.Range("E21").Value = 0.04070966823459377
$ws.Range("C22").Value = 31.20880547241887
$ws.Range("D22").Value = 0.08880547241886561
$ws.Range("E22").Value = 0.0078864119315379
$ws.Range("C23").Value = 31.71783149681082
$ws.Range("D23").Value = 0.4378314968108157
$ws.Range("E23").Value = 0.1916964195995993
$ws.Range("C24").Value = 31.67425022380877
$ws.Range("D24").Value = 0.2942502238087705
$ws.Range("E24").Value = 0.08658319421151155
$ws.Range("C25").Value = 31.80670018641877
$ws.Range("D25").Value = 0.2267001864187677
$ws.Range("E25").Value = 0.05139297452230401
$ws.Range("C26").Value = 31.57236529702934
$ws.Range("D26").Value = -0.07763470297066277
$ws.Range("E26").Value = 0.006027147105343035
$ws.Range("C27").Value = 31.80795770925392
$ws.Range("D27").Value = -0.07204229074607937
$ws.Range("E27").Value = 0.005190091655942634
$ws.Range("C28").Value = 32.40764967107331
$ws.Range("D28").Value = 0.127649671073307
$ws.Range("E28").Value = 0.01629443852512348
$ws.Range("C29").Value = 32.0779747708137
$ws.Range("D29").Value = -0.3720252291863062
$ws.Range("E29").Value = 0.1384027711511236
$ws.Range("C30").Value = 33.49615976197383
$ws.Range("D30").Value = 0.6461597619738271
$ws.Range("E30").Value = 0.4175224379940728
$ws.Range("C31").Value = 32.90623986695228
$ws.Range("D31").Value = 0.006239866952284956
$ws.Range("E31").Value = 0.00003893593958221794
$ws.Range("C32").Value = 33.22022286182911
$ws.Range("D32").Value = 0.120222861829113
$ws.Range("E32").Value = 0.014453536506382
$ws.Range("C33").Value = 33.03518120681197
$ws.Range("D33").Value = -0.3648187931880287
$ws.Range("E33").Value = 0.1330927518631696
$ws.Range("C34").Value = 33.75397210975283
$ws.Range("D34").Value = 0.05397210975282718
$ws.Range("E34").Value = 0.002912988631171223
$ws.Range("C35").Value = 34.66480062472593
$ws.Range("D35").Value = 0.5648006247259261
$ws.Range("E35").Value = 0.3189997456907964
$ws.Range("C36").Value = 34.41753274239829
$ws.Range("D36").Value = 0.01753274239828784
$ws.Range("E36").Value = 0.00030739705600472
$ws.Range("C37").Value = 34.92800287556167
$ws.Range("D37").Value = 0.02800287556166836
$ws.Range("E37").Value = 0.0007841610397222831
$ws.Range("C38").Value = 34.47907648458726
$ws.Range("D38").Value = -0.820923515412737
$ws.Range("E38").Value = 0.6739154181576062
$ws.Range("C39").Value = 34.76454817958464
$ws.Range("D39").Value = -0.9354518204153663
$ws.Range("E39").Value = 0.8750701083184228
$ws.Range("C40").Value = 35.6386607338206
$ws.Range("D40").Value = -0.6613392661793966
$ws.Range("E40").Value = 0.4373696249907029
$ws.Range("C41").Value = 35.81547822064613
$ws.Range("D41").Value = -0.9845217793538694
$ws.Range("E41").Value = 0.9692831340221092
$ws.Range("C42").Value = 36.8833122197145
$ws.Range("D42").Value = -0.4166877802854998
$ws.Range("E42").Value = 0.1736287062392569
$ws.Range("C43").Value = 37.77573340676048
$ws.Range("D43").Value = -0.1242665932395184
$ws.Range("E43").Value = 0.01544218619535592
$ws.Range("C44").Value = 38.52955548243027
$ws.Range("D44").Value = 0.0295554824302684
$ws.Range("E44").Value = 0.0008735265416859039
$ws.Range("C45").Value = 39.2027103573383
$ws.Range("D45").Value = 0.3027103573383059
$ws.Range("E45").Value = 0.09163356043988483
$ws.Range("C46").Value = 40.11827417529376
$ws.Range("D46").Value = 0.7182741752937645
$ws.Range("E46").Value = 0.5159177908939375
$ws.Range("C47").Value = 40.63573555162989
$ws.Range("D47").Value = 0.7357355516298938
$ws.Range("E47").Value = 0.5413068019321441
$ws.Range("C48").Value = 39.12693659218619
$ws.Range("D48").Value = -0.973063407813811
$ws.Range("E48").Value = 0.9468523956262269
$ws.Range("C49").Value = 40.86997228496865
$ws.Range("D49").Value = 0.2699722849686523
$ws.Range("E49").Value = 0.07288503465119522
$ws.Range("C50").Value = 40.65588761654695
$ws.Range("D50").Value = -0.2441123834530501
$ws.Range("E50").Value = 0.05959085575512896
$ws.Range("C51").Value = 40.29513997690747
$ws.Range("D51").Value = -0.9048600230925317
$ws.Range("E51").Value = 0.818771661391017

# Update TOTAL row (52) and MSE row (53)
$ws.Range("C52").Value = -2.080545998229663
$ws.Range("E52").Value = 10.74188383000859
$ws.Range("E53").Value = 0.2148376766001719
